$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country data table updated (new countries inserted into list, shifting some
# rows; numeric COVID-19 stats refreshed to the newer snapshot time).
$updates = @(
    @{ Row=4; A="Estados Unidos"; B=1617389; C=24666; D=381319; E=1139893; F=0; G=1241; H=96177 },
    @{ Row=6; A="Brasil"; B=296940; C=3583; D=116683; E=161026; F=0; G=337; H=19231 },
    @{ Row=11; A="Alemania"; B=179021; C=490; D=158000; E=12712; F=0; G=39; H=8309 },
    @{ Row=39; A="Colombia"; B=18330; C=643; D=4431; E=13247; F=0; G=22; H=652 },
    @{ Row=64; A="Oman"; B=6370; C=327; D=1821; E=4518; F=0; G=1; H=31 },
    @{ Row=68; A="Camerun"; B=4288; C=555; D=1808; E=2324; F=0; G=10; H=156 },
    @{ Row=69; A="Luxemburgo"; B=3980; C=9; D=3741; E=130; F=0; G=0; H=109 },
    @{ Row=70; A="Irak"; B=3877; C=153; D=2483; E=1254; F=0; G=6; H=140 },
    @{ Row=71; A="Azerbaiyan"; B=3749; C=118; D=2340; E=1365; F=0; G=1; H=44 },
    @{ Row=75; A="Guinea"; B=3067; C=204; D=1575; E=1474; F=0; G=0; H=18 },
    @{ Row=76; A="Tailandia"; B=3037; C=3; D=2897; E=84; F=0; G=0; H=56 },
    @{ Row=77; A="Uzbekistan"; B=2964; C=25; D=2407; E=544; F=0; G=0; H=13 },
    @{ Row=164; A="Guyana"; B=127; C=2; D=57; E=60; F=0; G=0; H=10 },
    @{ Row=198; A="Belice"; B=18; C=0; D=16; E=0; F=0; G=0; H=2 },
    @{ Row=199; A="Nueva Caledonia"; B=18; C=0; D=18; E=0; F=0; G=0; H=0 },
    @{ Row=200; A="Santa Lucia"; B=18; C=0; D=18; E=0; F=0; G=0; H=0 },
    @{ Row=209; A="Montserrat"; B=11; C=0; D=10; E=0; F=0; G=0; H=1 },
    @{ Row=210; A="Groenlandia"; B=11; C=0; D=11; E=0; F=0; G=0; H=0 },
    @{ Row=211; A="Seychelles"; B=11; C=0; D=11; E=0; F=0; G=0; H=0 },
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}

# Update the "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 23:35"
